$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text updates (masthead volume/number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 33   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/12/2026  Through  1/18/2026"

# --- Cells whose role (style) changes: count <-> blank("0") <-> pct("***.* ") ---
# Pattern: (re)apply correct format from a same-column donor with the target role,
# forcing a literal-text write via NumberFormat "@" when the target is the placeholder text,
# so Excel stores it as a string rather than re-parsing "0" as a number.
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("C17").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 2

$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 5

$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -80

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D16").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 2

$ws.Range("E16").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100

$ws.Range("J15").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = 2

$ws.Range("K15").Copy()
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("K22").Value = -50

$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("L16").Copy()
$ws.Range("L27").PasteSpecial(-4122)
$ws.Range("L27").Value = -100

# --- Plain value updates (same style/role, new number) ---
$ws.Range("G15").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = -25
$ws.Range("M16").Value = -70
$ws.Range("N16").Value = -93.478260869565
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = -55.555555555555
$ws.Range("L17").Value = -20
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -50
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = -45.454545454545
$ws.Range("M18").Value = -25
$ws.Range("N18").Value = -83.333333333333
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -24
$ws.Range("I19").Value = 11
$ws.Range("J19").Value = 18
$ws.Range("K19").Value = -38.888888888888
$ws.Range("L19").Value = -45
$ws.Range("M19").Value = -52.173913043478
$ws.Range("N19").Value = -42.105263157894
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 50
$ws.Range("L20").Value = -75
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -96.153846153846
$ws.Range("C21").Value = 8
$ws.Range("E21").Value = -42.857142857142
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 52
$ws.Range("H21").Value = -19.230769230769
$ws.Range("I21").Value = 26
$ws.Range("J21").Value = 37
$ws.Range("K21").Value = -29.729729729729
$ws.Range("L21").Value = -45.833333333333
$ws.Range("M21").Value = -42.222222222222
$ws.Range("N21").Value = -84.242424242424
$ws.Range("M22").Value = -50
$ws.Range("L23").Value = -50
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 42
$ws.Range("E24").Value = -40.476190476190
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -9.174311926605
$ws.Range("I24").Value = 61
$ws.Range("J24").Value = 75
$ws.Range("K24").Value = -18.666666666666
$ws.Range("L24").Value = 27.083333333333
$ws.Range("M24").Value = 48.780487804878
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -60.714285714285
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = -20.547945205479
$ws.Range("I25").Value = 38
$ws.Range("J25").Value = 47
$ws.Range("K25").Value = -19.148936170212
$ws.Range("L25").Value = 5.555555555555
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 9
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -70.967741935483
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = -78.947368421052
$ws.Range("L26").Value = -55.555555555555
$ws.Range("M26").Value = -63.636363636363
$ws.Range("G27").Value = 1
$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = 50
$ws.Range("L28").Value = -66.666666666666
$ws.Range("J41").Value = 78
$ws.Range("K41").Value = -68.924302788844
$ws.Range("L41").Value = -83.189655172413
$ws.Range("M41").Value = -90.982658959537
$ws.Range("N41").Value = -93.004484304932
$ws.Range("J43").Value = 119
$ws.Range("K43").Value = -65.406976744186
$ws.Range("L43").Value = -76.104417670682
$ws.Range("M43").Value = -85.069008782936
$ws.Range("N43").Value = -92.761557177615
$ws.Range("J44").Value = 407
$ws.Range("K44").Value = 10.298102981029
$ws.Range("L44").Value = -9.955752212389
$ws.Range("M44").Value = -20.196078431372
$ws.Range("N44").Value = -26
$ws.Range("J46").Value = 848
$ws.Range("K46").Value = -36.144578313253
$ws.Range("L46").Value = -56.624040920716
$ws.Range("M46").Value = -75.819788993441
$ws.Range("N46").Value = -83.607191184999

$excel.CutCopyMode = $false
